$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new "Scale" column (column D) to the position table.
# D9 is the header, D10:D64 hold a scale value of 1 for each data row.
$ws.Range("D9").Value = "Scale"

for ($r = 10; $r -le 64; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Update the active selection to reflect the newly added column and
# scroll the view back to the top-left of the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D11:D64").Select()
